# Refresh currentAveragePrice / Leve profit columns (H:N) with latest market data
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1173.9193
$ws.Range("I40").Value = 1100
$ws.Range("J40").Value = 1189.8628
$ws.Range("K40").Value = 1100
$ws.Range("L40").Value = 1189.8628
$ws.Range("M40").Value = -925
$ws.Range("N40").Value = -1539.8628
$ws.Range("H74").Value = 4340
$ws.Range("I74").Value = 3844
$ws.Range("K74").Value = 3844
$ws.Range("M74").Value = -2908
$ws.Range("H76").Value = 8237.16
$ws.Range("I76").Value = 12366.272
$ws.Range("K76").Value = 12366.272
$ws.Range("M76").Value = -12051.272
$ws.Range("H77").Value = 4340
$ws.Range("I77").Value = 3844
$ws.Range("K77").Value = 19220
$ws.Range("M77").Value = -14540
$ws.Range("H79").Value = 8237.16
$ws.Range("I79").Value = 12366.272
$ws.Range("K79").Value = 12366.272
$ws.Range("M79").Value = -11274.272
$ws.Range("H98").Value = 1582.15
$ws.Range("I98").Value = 599.6875
$ws.Range("J98").Value = 5512
$ws.Range("K98").Value = 599.6875
$ws.Range("L98").Value = 5512
$ws.Range("M98").Value = 898.3125
$ws.Range("N98").Value = -8508
$ws.Range("H122").Value = 1582.15
$ws.Range("I122").Value = 599.6875
$ws.Range("J122").Value = 5512
$ws.Range("K122").Value = 1799.0625
$ws.Range("L122").Value = 16536
$ws.Range("M122").Value = 650.9375
$ws.Range("N122").Value = -21436
$ws.Range("H125").Value = 2169.5925
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 2176.1155
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 19585.0395
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -24505.0395
$ws.Range("H129").Value = 1054.3
$ws.Range("I129").Value = 308.55554
$ws.Range("J129").Value = 1373.9048
$ws.Range("K129").Value = 925.66662
$ws.Range("L129").Value = 4121.7144
$ws.Range("M129").Value = 4074.33338
$ws.Range("N129").Value = -14121.7144
$ws.Range("H137").Value = 5221.3516
$ws.Range("I137").Value = 6526.8076
$ws.Range("J137").Value = 2135.7273
$ws.Range("K137").Value = 19580.4228
$ws.Range("L137").Value = 6407.1819
$ws.Range("M137").Value = -17030.4228
$ws.Range("N137").Value = -11507.1819

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 870848
$ws.Range("I32").Value = 971699
$ws.Range("K32").Value = 971699
$ws.Range("M32").Value = -971412
$ws.Range("H63").Value = 3428.4285
$ws.Range("I63").Value = 3499.8333
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 3499.8333
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2813.8333
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 3428.4285
$ws.Range("I66").Value = 3499.8333
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 17499.1665
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -14067.1665
$ws.Range("N66").Value = -21864
$ws.Range("H74").Value = 1637.2203
$ws.Range("I74").Value = 1021.2558
$ws.Range("J74").Value = 3292.625
$ws.Range("K74").Value = 1021.2558
$ws.Range("L74").Value = 3292.625
$ws.Range("M74").Value = -147.2558
$ws.Range("N74").Value = -5040.625
$ws.Range("H77").Value = 1637.2203
$ws.Range("I77").Value = 1021.2558
$ws.Range("J77").Value = 3292.625
$ws.Range("K77").Value = 5106.279
$ws.Range("L77").Value = 16463.125
$ws.Range("M77").Value = -738.2790000000005
$ws.Range("N77").Value = -25199.125
$ws.Range("H110").Value = 2516.3076
$ws.Range("I110").Value = 2540.8333
$ws.Range("J110").Value = 2222
$ws.Range("K110").Value = 2540.8333
$ws.Range("L110").Value = 2222
$ws.Range("M110").Value = -495.8332999999998
$ws.Range("N110").Value = -6312

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1132.6086
$ws.Range("I20").Value = 979.0625
$ws.Range("J20").Value = 1483.5714
$ws.Range("K20").Value = 979.0625
$ws.Range("L20").Value = 1483.5714
$ws.Range("M20").Value = -732.0625
$ws.Range("N20").Value = -1977.5714
$ws.Range("H80").Value = 126.44444
$ws.Range("I80").Value = 132.6
$ws.Range("J80").Value = 118.75
$ws.Range("K80").Value = 132.6
$ws.Range("L80").Value = 118.75
$ws.Range("M80").Value = 865.4
$ws.Range("N80").Value = -2114.75
$ws.Range("H83").Value = 126.44444
$ws.Range("I83").Value = 132.6
$ws.Range("J83").Value = 118.75
$ws.Range("K83").Value = 663
$ws.Range("L83").Value = 593.75
$ws.Range("M83").Value = 4329
$ws.Range("N83").Value = -10577.75
$ws.Range("H94").Value = 1432
$ws.Range("I94").Value = 1900
$ws.Range("J94").Value = 1244.8
$ws.Range("K94").Value = 1900
$ws.Range("L94").Value = 1244.8
$ws.Range("M94").Value = -1449
$ws.Range("N94").Value = -2146.8
$ws.Range("H104").Value = 29500
$ws.Range("J104").Value = 29500
$ws.Range("L104").Value = 29500
$ws.Range("N104").Value = -36488

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3271266.8
$ws.Range("I62").Value = 9261573
$ws.Range("J62").Value = 3827.2727
$ws.Range("K62").Value = 9261573
$ws.Range("L62").Value = 3827.2727
$ws.Range("M62").Value = -9260949
$ws.Range("N62").Value = -5075.2727
$ws.Range("H65").Value = 3271266.8
$ws.Range("I65").Value = 9261573
$ws.Range("J65").Value = 3827.2727
$ws.Range("K65").Value = 46307865
$ws.Range("L65").Value = 19136.3635
$ws.Range("M65").Value = -46304745
$ws.Range("N65").Value = -25376.3635
$ws.Range("H132").Value = 2396.9143
$ws.Range("I132").Value = 1172
$ws.Range("K132").Value = 3516
$ws.Range("M132").Value = -986
$ws.Range("H134").Value = 2221.3333
$ws.Range("I134").Value = 1327.6666
$ws.Range("J134").Value = 3412.889
$ws.Range("K134").Value = 3982.9998
$ws.Range("L134").Value = 10238.667
$ws.Range("M134").Value = -1447.9998
$ws.Range("N134").Value = -15308.667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2978.682
$ws.Range("I131").Value = 4794.4443
$ws.Range("J131").Value = 2511.7715
$ws.Range("K131").Value = 14383.3329
$ws.Range("L131").Value = 7535.314499999999
$ws.Range("M131").Value = -9343.332900000001
$ws.Range("N131").Value = -17615.3145

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5165.1724
$ws.Range("I80").Value = 6154.5
$ws.Range("J80").Value = 2966.6667
$ws.Range("K80").Value = 6154.5
$ws.Range("L80").Value = 2966.6667
$ws.Range("M80").Value = -5156.5
$ws.Range("N80").Value = -4962.6667
$ws.Range("H83").Value = 5165.1724
$ws.Range("I83").Value = 6154.5
$ws.Range("J83").Value = 2966.6667
$ws.Range("K83").Value = 30772.5
$ws.Range("L83").Value = 14833.3335
$ws.Range("M83").Value = -25780.5
$ws.Range("N83").Value = -24817.3335
$ws.Range("H95").Value = 24172
$ws.Range("J95").Value = 24172
$ws.Range("L95").Value = 24172
$ws.Range("N95").Value = -29664
$ws.Range("H102").Value = 4175.8945
$ws.Range("I102").Value = 2360.6667
$ws.Range("J102").Value = 7287.7144
$ws.Range("K102").Value = 2360.6667
$ws.Range("L102").Value = 7287.7144
$ws.Range("M102").Value = -738.6667000000002
$ws.Range("N102").Value = -10531.7144
$ws.Range("H132").Value = 3763.4634
$ws.Range("I132").Value = 3549.926
$ws.Range("J132").Value = 4175.2856
$ws.Range("K132").Value = 10649.778
$ws.Range("L132").Value = 12525.8568
$ws.Range("M132").Value = -8119.778
$ws.Range("N132").Value = -17585.8568

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1711.3529
$ws.Range("I16").Value = 1711.3529
$ws.Range("K16").Value = 1711.3529
$ws.Range("M16").Value = -1541.3529
$ws.Range("H22").Value = 717.8946999999999
$ws.Range("I22").Value = 621.5714
$ws.Range("J22").Value = 774.0833
$ws.Range("K22").Value = 621.5714
$ws.Range("L22").Value = 774.0833
$ws.Range("M22").Value = -326.5714
$ws.Range("N22").Value = -1364.0833
$ws.Range("H27").Value = 717.8946999999999
$ws.Range("I27").Value = 621.5714
$ws.Range("J27").Value = 774.0833
$ws.Range("K27").Value = 621.5714
$ws.Range("L27").Value = 774.0833
$ws.Range("M27").Value = -514.5714
$ws.Range("N27").Value = -988.0833
$ws.Range("H46").Value = 1078
$ws.Range("I46").Value = 1223.5
$ws.Range("J46").Value = 932.5
$ws.Range("K46").Value = 1223.5
$ws.Range("L46").Value = 932.5
$ws.Range("M46").Value = -1035.5
$ws.Range("N46").Value = -1308.5
$ws.Range("H55").Value = 237.9
$ws.Range("I55").Value = 208.77777
$ws.Range("J55").Value = 500
$ws.Range("K55").Value = 208.77777
$ws.Range("L55").Value = 500
$ws.Range("M55").Value = -35.77777
$ws.Range("N55").Value = -846
$ws.Range("H68").Value = 3233.5715
$ws.Range("I68").Value = 2333.6667
$ws.Range("J68").Value = 3908.5
$ws.Range("K68").Value = 2333.6667
$ws.Range("L68").Value = 3908.5
$ws.Range("M68").Value = -1584.6667
$ws.Range("N68").Value = -5406.5
$ws.Range("H71").Value = 3233.5715
$ws.Range("I71").Value = 2333.6667
$ws.Range("J71").Value = 3908.5
$ws.Range("K71").Value = 11668.3335
$ws.Range("L71").Value = 19542.5
$ws.Range("M71").Value = -7924.333500000001
$ws.Range("N71").Value = -27030.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1996.4117
$ws.Range("I81").Value = 921.9091
$ws.Range("J81").Value = 3966.3333
$ws.Range("K81").Value = 1843.8182
$ws.Range("L81").Value = 7932.6666
$ws.Range("M81").Value = -782.8181999999999
$ws.Range("N81").Value = -10054.6666
$ws.Range("H84").Value = 1996.4117
$ws.Range("I84").Value = 921.9091
$ws.Range("J84").Value = 3966.3333
$ws.Range("K84").Value = 9219.091
$ws.Range("L84").Value = 39663.333
$ws.Range("M84").Value = -3915.091
$ws.Range("N84").Value = -50271.333
$ws.Range("H107").Value = 397.58334
$ws.Range("I107").Value = 433
$ws.Range("J107").Value = 348
$ws.Range("K107").Value = 1299
$ws.Range("L107").Value = 1044
$ws.Range("M107").Value = 621
$ws.Range("N107").Value = -4884
$ws.Range("H113").Value = 345.25
$ws.Range("I113").Value = 401.6
$ws.Range("J113").Value = 251.33333
$ws.Range("K113").Value = 1204.8
$ws.Range("L113").Value = 753.99999
$ws.Range("M113").Value = 965.1999999999998
$ws.Range("N113").Value = -5093.99999
$ws.Range("H132").Value = 1592.2075
$ws.Range("I132").Value = 825.25714
$ws.Range("J132").Value = 3083.5
$ws.Range("K132").Value = 2475.77142
$ws.Range("L132").Value = 9250.5
$ws.Range("M132").Value = 54.22857999999997
$ws.Range("N132").Value = -14310.5
